$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- 1. Update the first three summary cells ---
$t.Cell(1, 1).Range.Text = "0M"
$t.Cell(2, 1).Range.Text = "0M"
$t.Cell(3, 1).Range.Text = "0M"

# --- 2. Insert 10 new rows right after row 3 (i.e. before the current row 4) ---
$newValues = @("120", "0.00002", "0.00006", "0.00004", "0.00001", "0.00004", "0.00004", "0.00004", "0.00481", "100.0")

# Inserting repeatedly "before" the same anchor row stacks the new rows above it
# in the order they were created, so row 4 ends up holding $newValues[0], etc.
$refRow = $t.Rows.Item(4)
foreach ($val in $newValues) {
    $t.Rows.Add($refRow) | Out-Null
}
for ($i = 0; $i -lt $newValues.Count; $i++) {
    $t.Cell(4 + $i, 1).Range.Text = $newValues[$i]
}

# --- 3. Collapse the final three multi-run / tab-separated cells into single values ---
$n = $t.Rows.Count
$t.Cell($n - 2, 1).Range.Text = "100"
$t.Cell($n - 1, 1).Range.Text = "0"
$t.Cell($n, 1).Range.Text = "111"
